# "Teste 1,2,3 passando com excel"
# Consolidate the three test sheets into a single sheet: fold the row that
# lived on "Teste3" into "Teste1" (row 4) then drop the now-redundant
# "Teste2" / "Teste3" sheets, leaving just "Teste1".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Teste1")

# Row 4 on Teste1 carries the data that used to live on Teste3 (username +
# password), paired with the matching firstname/lastname from Teste1/Teste2.
$ws1.Cells.Item(4, 1).Value = "Fabiana"
$ws1.Cells.Item(4, 2).Value = "Zimmer"
$ws1.Cells.Item(4, 3).Value = "zimmerf"
$ws1.Cells.Item(4, 4).Value = 123456

# Drop the other two sheets now that their data lives on Teste1.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Teste2").Delete() | Out-Null
$wb.Worksheets.Item("Teste3").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Match the final selection left behind on the remaining sheet.
$ws1.Select()
$ws1.Range("B7").Select() | Out-Null
